# Readme / minor changes:
#  - Insert a new "Clicking on SoapAura Logo" test row above the existing
#    row 5 (pushing every subsequent row down by one).
#  - Update the print area to cover the extra row.
#  - Move the active selection to the newly inserted B5 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above what is currently row 5. Excel shifts rows
# 5..20 down to 6..21 automatically.
$ws.Rows(5).Insert()

# Clone the formatting of row 4 (the other "header" style test-case row:
# plain style on A/B/D, wrap style on C) onto the freshly inserted row 5.
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)
$ws.Range("A5").Value = $null

# Fill in the new test case content (order matches first-use order so the
# shared-string table is appended to in the same sequence as the source
# workbook).
$ws.Range("C5").Value = "1) Clicking on Logo will link to home page"
$ws.Range("D5").Value = "1) Link user to home page"
$ws.Range("B5").Value = "Clicking on SoapAura Logo"

$ws.Rows(5).RowHeight = 43.8

# Move the selection onto the new cell.
$ws.Range("B5").Select()

# Extend the print area by one row to keep the new row inside it.
$ws.PageSetup.PrintArea = "`$A`$2:`$D`$20"
